$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.746.80"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.757.41"
$ws.Range("E3").Value = "  -3.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.23"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4249"
$ws.Range("E7").Value = "  -3.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3641"
$ws.Range("E8").Value = "  -1.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07543"
$ws.Range("E9").Value = "  -2.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.63"
$ws.Range("E10").Value = "  -4.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.095"
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.66"
$ws.Range("E13").Value = "  -6.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.043"
$ws.Range("E14").Value = "  -3.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.268"
$ws.Range("E15").Value = "  -3.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.786.36"
$ws.Range("E16").Value = "  -2.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.37"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001073"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06390"
$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.05"
$ws.Range("E21").Value = "  -2.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.906"
$ws.Range("E22").Value = "  -4.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.796.64"
$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("E24").Value = "  -4.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.102"
$ws.Range("E25").Value = "  +5.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.66"
$ws.Range("E26").Value = "  +3.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.32"
$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.980.91"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.138"
$ws.Range("E29").Value = "  -7.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.28"
$ws.Range("E30").Value = "  -2.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.115"
$ws.Range("E31").Value = "  -7.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.667"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.572"
$ws.Range("E33").Value = "  -5.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08908"
$ws.Range("E34").Value = "  -3.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.23"
$ws.Range("E35").Value = "  -6.59%  "

$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2107"
$ws.Range("E37").Value = "  -2.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06035"
$ws.Range("E38").Value = "  -2.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6334"
$ws.Range("E39").Value = "  -3.78%  "

$ws.Range("E40").Value = "  -4.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.186"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.908"
$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.395"
$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.702"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5869"
$ws.Range("E47").Value = "  -3.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.988"
$ws.Range("E48").Value = "  -2.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.32"
$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.176"
$ws.Range("E50").Value = "  +1.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06829"
$ws.Range("E51").Value = "  -2.16%  "
